{"js": "// Update the benchmark results table:\n// - rows 1-3 (first column of results) become \"0M\" (not yet measured)\n// - row 4 (\"105\") becomes \"307\"\n// - the following summary-stat rows get refreshed numbers\n// - the last three rows collapse their multi-value (tab separated) runs\n//   down to the single \"total\" figure that used to live in the first rows\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, newValue) pairs, applied to column 0 of the single-column table.\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"307\"],\n  [4, \"0.00002\"],\n  [5, \"0.00056\"],\n  [6, \"0.00017\"],\n  [8, \"0.00030\"],\n  [9, \"0.00037\"],\n  [10, \"0.00044\"],\n  [11, \"0.06615\"],\n  [43, \"99.72\"],\n  [44, \"0.07\"],\n  [45, \"23\"],\n];\n\nfor (const [rowIndex, value] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = value;\n}\n\nawait context.sync();\n", "ps1": "# Update the benchmark results table:\n# - rows 1-3 (first column of results) become \"0M\" (not yet measured)\n# - row 4 (\"105\") becomes \"307\"\n# - the following summary-stat rows get refreshed numbers\n# - the last three rows collapse their multi-value (tab separated) runs\n#   down to the single \"total\" figure that used to live in the first rows\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"307\"\n$t.Cell(5, 1).Range.Text = \"0.00002\"\n$t.Cell(6, 1).Range.Text = \"0.00056\"\n$t.Cell(7, 1).Range.Text = \"0.00017\"\n$t.Cell(9, 1).Range.Text = \"0.00030\"\n$t.Cell(10, 1).Range.Text = \"0.00037\"\n$t.Cell(11, 1).Range.Text = \"0.00044\"\n$t.Cell(12, 1).Range.Text = \"0.06615\"\n$t.Cell(44, 1).Range.Text = \"99.72\"\n$t.Cell(45, 1).Range.Text = \"0.07\"\n$t.Cell(46, 1).Range.Text = \"23\"\n"}
